$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (both columns to 15.42578125)
$ws.Columns.Item(1).ColumnWidth = 15.42578125
$ws.Columns.Item(2).ColumnWidth = 15.42578125

# Update cell values
$ws.Cells.Item(1, 1).Value = -0.09476902780640728
$ws.Cells.Item(1, 2).Value = 0.094275394532772339
$ws.Cells.Item(2, 1).Value = -0.10737733275275918
$ws.Cells.Item(2, 2).Value = 0.10591829089725557
$ws.Cells.Item(3, 1).Value = -0.11559457599656575
$ws.Cells.Item(3, 2).Value = 0.11515636157261966
$ws.Cells.Item(4, 1).Value = -0.10715636160349007
$ws.Cells.Item(4, 2).Value = 0.106766347044017
$ws.Cells.Item(5, 1).Value = -0.103766347060537
$ws.Cells.Item(5, 2).Value = 0.10244854774759737
$ws.Cells.Item(6, 1).Value = -0.0031717058906064466
$ws.Cells.Item(6, 2).Value = 0.0031320637154284015
$ws.Cells.Item(7, 1).Value = 0.0068679362412327194
$ws.Cells.Item(7, 2).Value = -0.0068703572627875076
$ws.Cells.Item(8, 1).Value = -0.031943374029884897
$ws.Cells.Item(8, 2).Value = 0.031658228601652461
$ws.Cells.Item(9, 1).Value = -0.029658228620223159
$ws.Cells.Item(9, 2).Value = 0.029418653062637468
$ws.Cells.Item(10, 1).Value = -0.0274186530826519
$ws.Cells.Item(10, 2).Value = 0.027402096268350107
$ws.Cells.Item(11, 1).Value = -0.0244020962917606
$ws.Cells.Item(11, 2).Value = 0.024374867639466835
$ws.Cells.Item(12, 1).Value = -0.020874867664927077
$ws.Cells.Item(12, 2).Value = 0.020675721527586344
$ws.Cells.Item(13, 1).Value = -0.017175721554456125
$ws.Cells.Item(13, 2).Value = 0.017084857564388045
$ws.Cells.Item(14, 1).Value = -0.0090848576062434461
$ws.Cells.Item(14, 2).Value = 0.0090550196316838694
$ws.Cells.Item(15, 1).Value = -0.0080550196509809879
$ws.Cells.Item(15, 2).Value = 0.0080357355505631389
$ws.Cells.Item(16, 1).Value = -0.0060357355733562379
$ws.Cells.Item(16, 2).Value = 0.0060037478987751136
$ws.Cells.Item(17, 1).Value = -0.0040037479218302252
$ws.Cells.Item(17, 2).Value = 0.0039999999703317357
$ws.Cells.Item(18, 1).Value = 0.019101938200147828
$ws.Cells.Item(18, 2).Value = -0.019211077464245108
$ws.Cells.Item(19, 1).Value = -0.012092114837048751
$ws.Cells.Item(19, 2).Value = 0.012017084742538486
$ws.Cells.Item(20, 1).Value = -0.0080170847566716219
$ws.Cells.Item(20, 2).Value = 0.0080057293869160162
$ws.Cells.Item(21, 1).Value = -0.0040057294011912603
$ws.Cells.Item(21, 2).Value = 0.0039999999856084045
$ws.Cells.Item(22, 1).Value = 0.013676285083427331
$ws.Cells.Item(22, 2).Value = -0.013964137296337498
$ws.Cells.Item(23, 1).Value = 0.018964137277481363
$ws.Cells.Item(23, 2).Value = -0.01953999862276401
$ws.Cells.Item(24, 1).Value = -0.02009943012395432
$ws.Cells.Item(24, 2).Value = 0.019999999932340806
$ws.Cells.Item(25, 1).Value = -0.09727684187458685
$ws.Cells.Item(25, 2).Value = 0.097150172539102186
$ws.Cells.Item(26, 1).Value = -0.094650172559241241
$ws.Cells.Item(26, 2).Value = 0.094487408752966573
$ws.Cells.Item(27, 1).Value = -0.091987408774524937
$ws.Cells.Item(27, 2).Value = 0.091025994951676736
$ws.Cells.Item(28, 1).Value = -0.089025994977252054
$ws.Cells.Item(28, 2).Value = 0.088367671991705699
$ws.Cells.Item(29, 1).Value = -0.081367672037636396
$ws.Cells.Item(29, 2).Value = 0.081176600465624027
$ws.Cells.Item(30, 1).Value = -0.021176600686123148
$ws.Cells.Item(30, 2).Value = 0.021024744088222747
$ws.Cells.Item(31, 1).Value = -0.014024744137623557
$ws.Cells.Item(31, 2).Value = 0.014001512907075409
$ws.Cells.Item(32, 1).Value = -0.0040015129664272564
$ws.Cells.Item(32, 2).Value = 0.0039999999601025849
